$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.044578486661111
$ws.Range("D2").Value = 1.050286781194331
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.05803654146201
$ws.Range("I2").Value = 1.04515051847246
$ws.Range("J2").Value = 1.049642820311586
$ws.Range("K2").Value = 1.053041479020159
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.060769873490923
$ws.Range("N2").Value = 1.02040366397046

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045554242169184
$ws.Range("D3").Value = 1.051061122918184
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.058996407795321
$ws.Range("I3").Value = 1.045441573898956
$ws.Range("J3").Value = 1.050265758189885
$ws.Range("K3").Value = 1.05362833502183
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.06154332365306
$ws.Range("N3").Value = 1.02061253897232

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.04618565559105
$ws.Range("D4").Value = 1.051562195828748
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.059617964197831
$ws.Range("I4").Value = 1.045628581835255
$ws.Range("J4").Value = 1.050668256196563
$ws.Range("K4").Value = 1.054007432060348
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.062043633019576
$ws.Range("N4").Value = 1.020747442295809

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.046451109547178
$ws.Range("D5").Value = 1.051772851209962
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.059879375418846
$ws.Range("I5").Value = 1.045706882512996
$ws.Range("J5").Value = 1.050837325841185
$ws.Range("K5").Value = 1.054166651052202
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.062253922896135
$ws.Range("N5").Value = 1.020804094888315

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.046495680875457
$ws.Range("D6").Value = 1.05180822140176
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.059923273873336
$ws.Range("I6").Value = 1.04572001091385
$ws.Range("J6").Value = 1.050865705121084
$ws.Range("K6").Value = 1.054193375606826
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.062289229112998
$ws.Range("N6").Value = 1.020813603537309

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.046189202569921
$ws.Range("D7").Value = 1.051565010600437
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.059621456760553
$ws.Range("I7").Value = 1.045629629340635
$ws.Range("J7").Value = 1.050670515866315
$ws.Range("K7").Value = 1.054009560155902
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.062046443082052
$ws.Range("N7").Value = 1.020748199529601

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.04490824040447
$ws.Range("D8").Value = 1.050548468474558
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.05836083726304
$ws.Range("I8").Value = 1.045249156000512
$ws.Range("J8").Value = 1.049853465740451
$ws.Range("K8").Value = 1.053239941469285
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.061031298212719
$ws.Range("N8").Value = 1.020474306442118

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.042651305091643
$ws.Range("D9").Value = 1.0487574016253
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.056143017386474
$ws.Range("I9").Value = 1.044568587552496
$ws.Range("J9").Value = 1.048409268920933
$ws.Range("K9").Value = 1.051878917124295
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.059241256443535
$ws.Range("N9").Value = 1.019989746912642

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.041146897546715
$ws.Range("D10").Value = 1.047563552870543
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.054666911092268
$ws.Range("I10").Value = 1.044108089731261
$ws.Range("J10").Value = 1.047443512716643
$ws.Range("K10").Value = 1.050968338687629
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.058047114465041
$ws.Range("N10").Value = 1.019665427359175

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.040495528244239
$ws.Range("D11").Value = 1.047046660628641
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.0540283298747
$ws.Range("I11").Value = 1.043907084068438
$ws.Range("J11").Value = 1.047024634537296
$ws.Range("K11").Value = 1.050573289727384
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.057529861358864
$ws.Range("N11").Value = 1.019524692705537

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.040253588327366
$ws.Range("D12").Value = 1.046854672499619
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.053791220661578
$ws.Range("I12").Value = 1.043832180472088
$ws.Range("J12").Value = 1.046868939921693
$ws.Range("K12").Value = 1.050426436931948
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.057337703856139
$ws.Range("N12").Value = 1.019472372402872

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.040305484910565
$ws.Range("D13").Value = 1.046895854199638
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.053842077385408
$ws.Range("I13").Value = 1.043848258455189
$ws.Range("J13").Value = 1.046902341671611
$ws.Range("K13").Value = 1.050457942506255
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.05737892348934
$ws.Range("N13").Value = 1.019483597326827

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.040475529257733
$ws.Range("D14").Value = 1.047030790648466
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.054008728547803
$ws.Range("I14").Value = 1.04390089743131
$ws.Range("J14").Value = 1.047011766903119
$ws.Range("K14").Value = 1.050561153152435
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.057513978100087
$ws.Range("N14").Value = 1.019520368813577

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.040580300213218
$ws.Range("D15").Value = 1.04711393062094
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.054111419550419
$ws.Range("I15").Value = 1.043933298096964
$ws.Range("J15").Value = 1.047079173551858
$ws.Range("K15").Value = 1.050624729544846
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.057597186188376
$ws.Range("N15").Value = 1.019543018958801

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041190127839293
$ws.Range("D16").Value = 1.047597858475843
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.054709303963292
$ws.Range("I16").Value = 1.044121395981656
$ws.Range("J16").Value = 1.04747129758648
$ws.Range("K16").Value = 1.050994540766828
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.058081439075328
$ws.Range("N16").Value = 1.019674761106515

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.041572669930288
$ws.Range("D17").Value = 1.047901428181728
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.055084497562981
$ws.Range("I17").Value = 1.044238954769465
$ws.Range("J17").Value = 1.047717079719538
$ws.Range("K17").Value = 1.051226309812978
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.058385149908238
$ws.Range("N17").Value = 1.01975731872888

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.041795804978437
$ws.Range("D18").Value = 1.048078500231104
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.055303397798316
$ws.Range("I18").Value = 1.044307369742579
$ws.Range("J18").Value = 1.047860372769355
$ws.Range("K18").Value = 1.051361423171599
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.058562281719813
$ws.Range("N18").Value = 1.019805444026523

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.041871889055725
$ws.Range("D19").Value = 1.048138878042713
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.055378046623564
$ws.Range("I19").Value = 1.044330671167214
$ws.Range("J19").Value = 1.047909220531338
$ws.Range("K19").Value = 1.051407480834366
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.058622676088977
$ws.Range("N19").Value = 1.019821848555834

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.041531626294168
$ws.Range("D20").Value = 1.047868857504507
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.055044237026433
$ws.Range("I20").Value = 1.044226357862428
$ws.Range("J20").Value = 1.047690716602199
$ws.Range("K20").Value = 1.051201450797147
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.058352566423755
$ws.Range("N20").Value = 1.019748464092647

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.040425455222547
$ws.Range("D21").Value = 1.046991054981659
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.053959651490959
$ws.Range("I21").Value = 1.043885403219156
$ws.Range("J21").Value = 1.0469795467838
$ws.Range("K21").Value = 1.05053076333288
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.057474208614431
$ws.Range("N21").Value = 1.019509541772922

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.039730005928187
$ws.Range("D22").Value = 1.046439196233755
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.053278240678506
$ws.Range("I22").Value = 1.043669636470489
$ws.Range("J22").Value = 1.046531801085277
$ws.Range("K22").Value = 1.050108416071908
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.056921796218055
$ws.Range("N22").Value = 1.019359060497784

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040098672462166
$ws.Range("D23").Value = 1.046731741930552
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.053639420736252
$ws.Range("I23").Value = 1.043784150638811
$ws.Range("J23").Value = 1.046769216744758
$ws.Range("K23").Value = 1.0503323726747
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.05721465486412
$ws.Range("N23").Value = 1.019438858170575

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.041550172131968
$ws.Range("D24").Value = 1.047883574776216
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.055062428856402
$ws.Range("I24").Value = 1.044232050341644
$ws.Range("J24").Value = 1.047702629168714
$ws.Range("K24").Value = 1.05121268374378
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.058367289552169
$ws.Range("N24").Value = 1.01975246521174

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.043234740751122
$ws.Range("D25").Value = 1.049220404801717
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.056715950643663
$ws.Range("I25").Value = 1.044745728452847
$ws.Range("J25").Value = 1.048783152303912
$ws.Range("K25").Value = 1.052231346378805
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.059704165989711
$ws.Range("N25").Value = 1.020115243799425
